# Apply edits described by the diff: add column S (copy/derive from column R)
# for rows 10-18, 31-35, 37-43, and update row 32 (B:S) with new "fit time" values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column S, rows 10-18
$ws.Range("S10").Value = 0
$ws.Range("S11").Value = 1
$ws.Range("S12").Value = 4.641196107526275
$ws.Range("S13").Value = 2.978216895447799
$ws.Range("S14").Value = -0.000001070866799108457
$ws.Range("S15").Value = -0.000001070866799108457
$ws.Range("S16").Value = 1.354007388198719
$ws.Range("S17").Value = 2.488138500418057
$ws.Range("S18").Value = 4.360125204043042

# New values for column S, rows 31-35
$ws.Range("S31").Value = 314

# Row 32 "fit time" -- updated values for B:R plus new S
$ws.Range("B32").Value = 13.72654837499999
$ws.Range("C32").Value = 32.32895194
$ws.Range("D32").Value = 16.94519662900001
$ws.Range("E32").Value = 26.915427016
$ws.Range("F32").Value = 16.12169647799999
$ws.Range("G32").Value = 59.77828406200001
$ws.Range("H32").Value = 13.21787197099999
$ws.Range("I32").Value = 26.84762463300001
$ws.Range("J32").Value = 16.69609554800002
$ws.Range("K32").Value = 28.873552297
$ws.Range("L32").Value = 24.61528613199999
$ws.Range("M32").Value = 70.861748111
$ws.Range("N32").Value = 13.65448766600002
$ws.Range("O32").Value = 34.76529264699997
$ws.Range("P32").Value = 18.88788201799997
$ws.Range("Q32").Value = 25.72435966300003
$ws.Range("R32").Value = 16.927842951
$ws.Range("S32").Value = 51.64329313899998

$ws.Range("S33").Value = 7.360955662482409
$ws.Range("S34").Value = 3
$ws.Range("S35").Value = 2.453651887494136

# New values for column S, rows 37-43
$ws.Range("S37").Value = 505
$ws.Range("S38").Value = 247
$ws.Range("S39").Value = 258
$ws.Range("S40").Value = 0.7364341085271318
$ws.Range("S41").Value = 0.0310391363022942
$ws.Range("S42").Value = 2.490169574385275
$ws.Range("S43").Value = 0.01761846232865465

# The newly written S cells pick up the column-level style (s="2") by
# default, but the source workbook leaves these particular cells with no
# explicit style (matching their B:R neighbours on the same rows). Reset
# them to the workbook's default "Normal" style so the XML stays consistent
# with the rest of the sheet.
$sCells = "S10,S11,S12,S13,S14,S15,S16,S17,S18,S31,S32,S33,S34,S35,S37,S38,S39,S40,S41,S42,S43"
foreach ($addr in $sCells.Split(",")) {
    $ws.Range($addr).Style = "Normal"
}
